# example.xlsx: rename plasma_gas "air" -> "air_11" and ic_db_name
# "db_example" -> "db_example.h5" across all case rows (comment from the
# commit: "slight fix to write excel + retrieve_helper moved" — the values
# written into the example sheet are updated to match the new helper
# defaults / file extensions).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# plasma_gas (column F) cells that currently read "air" become "air_11"
$airRows = @(3, 4, 6, 7, 9, 11)
foreach ($r in $airRows) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq "air") {
        $cell.Value = "air_11"
    }
}

# ic_db_name (column G) cells that currently read "db_example" become
# "db_example.h5"
$dbRows = @(3, 4, 5, 8, 11)
foreach ($r in $dbRows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "db_example") {
        $cell.Value = "db_example.h5"
    }
}
